# [Fix] Proceedings filters and download file
#
# The banner title in A1 of the "Expedientes" sheet has two date
# placeholders that used to repeat the same {{Fecha}} token for both the
# start and end of the range. Split them into distinct
# {{FechaInicial}} / {{FechaFinal}} placeholders so each can be bound to
# its own value when the report is rendered.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Expedientes")

$cell = $ws.Range("A1")

# The cell holds rich text (3 runs with different fonts/sizes). Replace
# just the date-range portion of the text, then re-apply each run's
# original formatting so the banner keeps its look
# (bold 14pt title / blank 12pt spacer / bold 12pt subtitle).
$full = $cell.Characters()
$oldText = $full.Text
if ($oldText.Contains("DEL {{Fecha}} AL {{Fecha}}")) {
    $newText = $oldText.Replace("DEL {{Fecha}} AL {{Fecha}}", "DEL {{FechaInicial}} AL {{FechaFinal}}")
    $full.Text = $newText
}

$titleText = "Laboratorio Alfonso Ramos S.A. de C.V. MONTERREY`n{{Direccion}}`n{{Sucursal}}`n"
$spacerText = "`n"
$subtitleText = "Listado de {{Titulo}}`nDEL {{FechaInicial}} AL {{FechaFinal}}"

$titleLen = $titleText.Length
$spacerLen = $spacerText.Length
$subtitleLen = $subtitleText.Length

$titleRun = $cell.Characters(1, $titleLen)
$titleRun.Font.Name = "Calibri"
$titleRun.Font.Bold = $true
$titleRun.Font.Size = 14

$spacerRun = $cell.Characters($titleLen + 1, $spacerLen)
$spacerRun.Font.Name = "Calibri"
$spacerRun.Font.Bold = $false
$spacerRun.Font.Size = 12

$subtitleRun = $cell.Characters($titleLen + $spacerLen + 1, $subtitleLen)
$subtitleRun.Font.Name = "Calibri"
$subtitleRun.Font.Bold = $true
$subtitleRun.Font.Size = 12

# Move the saved selection from F3 to the header band (A1:H1), matching
# the refreshed view state saved with the fix.
$ws.Range("A1:H1").Select()
